$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. "Semestre ideal" value: EA-7 -> EA-9
$ws.Range("B9").Value = "EA-9"
$ws.Range("C9").Value = "EA-9"

# 2. New "Requisito fraco" strings, replacing the old 3 "Requisito" rows (24-26)
# and continuing into new rows 27-51.
$newReq = @(
    "LOB1003 -  Cálculo I  (Requisito fraco)`n",
    "LOB1004 -  Cálculo II  (Requisito fraco)`n",
    "LOB1006 -  Cálculo IV  (Requisito fraco)`n",
    "LOB1011 -  Eletricidade Aplicada  (Requisito fraco)`n",
    "LOB1012 -  Estatística  (Requisito fraco)`n",
    "LOB1018 -  Física I  (Requisito fraco)`n",
    "LOB1019 -  Física II  (Requisito fraco)`n",
    "LOB1021 -  Física IV  (Requisito fraco)`n",
    "LOB1024 -  Mecânica  (Requisito fraco)`n",
    "LOB1036 -  Geometria Analítica  (Requisito fraco)`n",
    "LOB1037 -  Àlgebra Linear  (Requisito fraco)`n",
    "LOB1038 -  Física Experimental I  (Requisito fraco)`n",
    "LOB1039 -  Física Experimental III  (Requisito fraco)`n",
    "LOB1041 -  Física Experimental II  (Requisito fraco)`n",
    "LOB1042 -  Física Experimental IV  (Requisito fraco)`n",
    "LOB1045 -  Leitura e Produção de Textos Acadêmicos  (Requisito fraco)`n",
    "LOB1052 -  Cálculo III  (Requisito fraco)`n",
    "LOB1053 -  Física III  (Requisito fraco)`n",
    "LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito fraco)`n",
    "LOB1232 -  Licenciamento Ambiental  (Requisito fraco)`n",
    "LOB1257 -  Sistema de Abastecimento e Tratamento de Água  (Requisito fraco)`n",
    "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito fraco)`n",
    "LOM3081 -  Introdução à Mecânica dos Sólidos  (Requisito fraco)`n",
    "LOQ4095 -  Química Geral Experimental  (Requisito fraco)`n",
    "LOQ4097 -  Fundamentos de Química para Engenharia I (Requisito fraco)`n",
    "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito fraco)`n",
    "LOQ4233 -  Gestão de Negócios  (Requisito fraco)`n",
    "LOQ4247 -  Desenho Assistido por Computador  (Requisito fraco)`n",
)

$startRow = 24
for ($i = 0; $i -lt $newReq.Count; $i++) {
    $row = $startRow + $i
    if ($row -gt 26) {
        # New rows: copy formatting (styles) from row 26 template, then set value + height
        $ws.Range("B26:C26").Copy() | Out-Null
        $ws.Range("B" + $row + ":C" + $row).PasteSpecial(-4122) | Out-Null
    }
    $ws.Range("B" + $row).Value = $newReq[$i]
    $ws.Range("C" + $row).Value = $newReq[$i]
    $ws.Rows.Item($row).RowHeight = 30
}

$excel.CutCopyMode = $false
